$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "313.14"
Set-TextValue "E2" "1.50%"
Set-TextValue "D3" "39.18"
Set-TextValue "E3" "0.49%"
Set-TextValue "D4" "5.125"
Set-TextValue "E4" "0.02%"
Set-TextValue "D5" "0.08130"
Set-TextValue "E5" "0.14%"
Set-TextValue "D6" "4.493"
Set-TextValue "E6" "6.05%"
Set-TextValue "D7" "1.962"
Set-TextValue "E7" "0.86%"
Set-TextValue "D8" "8.286"
Set-TextValue "E8" "1.71%"
Set-TextValue "D9" "0.9388"
Set-TextValue "E9" "1.30%"
Set-TextValue "D10" "0.1325"
Set-TextValue "E10" "-5.82%"
Set-TextValue "D11" "0.1964"
Set-TextValue "E11" "1.78%"
Set-TextValue "D12" "0.09076"
Set-TextValue "E12" "0.44%"
Set-TextValue "D13" "0.03486"
Set-TextValue "E13" "-0.27%"
Set-TextValue "D14" "0.09710"
Set-TextValue "E14" "-1.16%"
Set-TextValue "D15" "0.001407"
Set-TextValue "E15" "0.91%"
Set-TextValue "E16" "3.77%"
Set-TextValue "D17" "3.556"
Set-TextValue "E17" "-8.92%"
Set-TextValue "D18" "3.170"
Set-TextValue "E18" "-5.57%"
Set-TextValue "D19" "0.3466"
Set-TextValue "E19" "0.33%"
Set-TextValue "E20" "-3.09%"
Set-TextValue "E21" "5.64%"
Set-TextValue "D23" "0.04373"
Set-TextValue "E24" "1.19%"
Set-TextValue "D25" "0.004728"
Set-TextValue "E25" "-1.47%"
Set-TextValue "E26" "199.37%"
Set-TextValue "E27" "-7.62%"
Set-TextValue "D39" "0.02212"
Set-TextValue "E39" "6.18%"
Set-TextValue "D40" "0.05225"
Set-TextValue "E40" "2.11%"
Set-TextValue "D41" "0.007608"
Set-TextValue "E41" "2.35%"
Set-TextValue "D42" "0.01034"
Set-TextValue "E42" "5.63%"
Set-TextValue "D43" "0.1394"
Set-TextValue "E43" "2.18%"
Set-TextValue "E44" "-1.36%"
Set-TextValue "D45" "0.009098"
Set-TextValue "E45" "-1.16%"
Set-TextValue "D46" "0.00006696"
Set-TextValue "E46" "4.71%"
Set-TextValue "E47" "0.06%"
Set-TextValue "D48" "0.003015"
Set-TextValue "E48" "17.09%"
Set-TextValue "D50" "0.00002103"
Set-TextValue "E50" "0.06%"
Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "0.06%"

Write-Host "Applied all updates"
